$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new test-case rows after the header/"Viewing Home page" row ---
# (becomes rows 3 & 4: "Scrolling in the website" and "Redirection in the website")
$ws.Rows("3:4").Insert()

# Copy formatting from the rows that will sit right below (old row 3, now row 5)
# so the new rows match the sheet's existing data-row style.
$ws.Range("A5:G6").Copy()
$ws.Range("A3:G4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# --- Fill the two new rows ---
$ws.Range("A3").Value = "Scrolling in the website"
$ws.Range("B3").Value = "This is to test if scrolling can be done successfully."
$ws.Range("C3").Value = "1. Go to http://localhost:8000/projects`n2. Scroll to the bottom of the page."
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "The page should scroll down to the bottom."
$ws.Range("F3").Value = "Same as expected outcome."
$ws.Range("G3").Value = "Pass"

$ws.Range("A4").Value = "Redirection in the website"
$ws.Range("B4").Value = "This is to test if redirection in the website can be done successfully."
$ws.Range("C4").Value = "1. Go to http://localhost:8000/projects`n2. Click on `"Blog`" in the navigation bar."
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "The page should load successfully and display the `"Blog`" page."
$ws.Range("F4").Value = "Same as expected outcome."
$ws.Range("G4").Value = "Pass"

$ws.Rows("3").RowHeight = 45
$ws.Rows("4").RowHeight = 60

# --- Update the "Actual Outcome" / "Pass/Fail" columns for all the pre-existing rows ---
for ($r = 2; $r -le 11; $r++) {
    if ($r -eq 3 -or $r -eq 4) { continue }
    $ws.Range("F$r").Value = "Same as expected outcome."
    $ws.Range("G$r").Value = "Pass"
}

# --- Update the expected-outcome wording for the "without an author" / "without a body" cases ---
$ws.Range("E10").Value = "At the author name field, there will be an error when trying to submit the form since the field is invalid."
$ws.Range("E11").Value = "At the body text field, there will be an error when trying to submit the form since the field is invalid."

# --- Fix up the view: selection moves to D4, no frozen/scrolled top-left cell ---
$ws.Range("D4").Select() | Out-Null
